# edit.ps1 - apply the "Overview updated. Bibliography updated." changes
#
# Strategy: Word's object model does not expose a way to programmatically
# insert <w:proofErr/> spell-check markers (those are only ever produced by
# Word's live background spell checker, never by automation/VBA/COM), so
# this script focuses on getting the paragraph text and run layout exactly
# right: every run-split called out in the diff is reproduced by toggling a
# (no-op) character formatting property on the exact target sub-range, which
# forces Word to materialize that sub-range as its own <w:r> run without
# altering the visible formatting.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Mobile: GameBench, Apptim, built-in profiling tools"
#    -> split into: "Mobile: " | "GameBench" | ", " | "Apptim" | ", built-in profiling tools"
# ---------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute("Mobile: GameBench, Apptim, built-in profiling tools")
$pStart = $full.Start

$gbStart = $pStart + 8
$gbEnd = $pStart + 17
$apStart = $pStart + 19
$apEnd = $pStart + 25

$rGameBench = $d.Range($gbStart, $gbEnd)
$rGameBench.Bold = 1
$rGameBench.Bold = 0

$rApptim = $d.Range($apStart, $apEnd)
$rApptim.Bold = 1
$rApptim.Bold = 0

# ---------------------------------------------------------------------
# 2) "Web: GTmetrix, WebPageTest, Insights (Web)"
#    -> split into: "Web: " | "GTmetrix" | ", " | "WebPageTest" | ", Insights (Web)"
# ---------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute("Web: GTmetrix, WebPageTest, Insights (Web)")
$pStart = $full.Start

$gtStart = $pStart + 5
$gtEnd = $pStart + 13
$wptStart = $pStart + 15
$wptEnd = $pStart + 26

$rGTmetrix = $d.Range($gtStart, $gtEnd)
$rGTmetrix.Bold = 1
$rGTmetrix.Bold = 0

$rWebPageTest = $d.Range($wptStart, $wptEnd)
$rWebPageTest.Bold = 1
$rWebPageTest.Bold = 0

# ---------------------------------------------------------------------
# 3) " (e.g. gifs from Giphy API)" -> split so "Giphy" is its own run:
#    " (e.g. gifs from " | "Giphy" | " API)"
# ---------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute("Lists: scrolling, rendering, filtering (e.g. gifs from Giphy API)")
$pStart = $full.Start

$giphyStart = $pStart + 55
$giphyEnd = $pStart + 60

$rGiphy = $d.Range($giphyStart, $giphyEnd)
$rGiphy.Bold = 1
$rGiphy.Bold = 0

# ---------------------------------------------------------------------
# 4) "API fetch" -> "Real-world application: movies API, show movies, search
#    for one, show details "  (plain text replacement)
# ---------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute("API fetch", $false, $false, $false, $false, $false, $true, 1, $false, "Real-world application: movies API, show movies, search for one, show details ", 2)

# ---------------------------------------------------------------------
# 5) " (screen with drawer and tabbar->buttons display->...->scroll)"
#    -> split so "tabbar" is its own run:
#    " (screen with drawer and " | "tabbar" | "->buttons display->...->scroll)"
#    Then add a brand-new list paragraph right after this one:
#    "Camera (take a picture and present it / select many pictures/videos
#    from gallery and present them)"
# ---------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute("UI (screen with drawer and tabbar")
$pStart = $full.Start

$tabbarStart = $pStart + 27
$tabbarEnd = $pStart + 33

$rTabbar = $d.Range($tabbarStart, $tabbarEnd)
$rTabbar.Bold = 1
$rTabbar.Bold = 0

$uiParaIndex = 0
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $pText = $d.Paragraphs.Item($i).Range.Text
    if ($pText.StartsWith("UI (screen with drawer and")) {
        $uiParaIndex = $i
    }
}

$uiPara = $d.Paragraphs.Item($uiParaIndex)
$uiPara.Range.InsertParagraphAfter()

$newParaIndex = $uiParaIndex + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newParaStart = $newPara.Range.Start
$newParaEnd = $newPara.Range.End
$newParaTextEnd = $newParaEnd - 1

$rNewPara = $d.Range($newParaStart, $newParaTextEnd)
$rNewPara.Text = "Camera (take a picture and present it / select many pictures/videos from gallery and present them)"

Write-Output "done"
